# Scheduled runner update: refresh market-board derived columns
# (currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ /
#  LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ) across sheets.
#
# Columns written, in order, for every row below: H, I, J, K, L, M, N.
# A $null entry means "leave that cell untouched" (matches rows that have
# no LeveProfitNQ value).

$wb = $excel.ActiveWorkbook

function Set-RowValues {
    param($ws, $row, $values)
    $col = 8  # column H
    foreach ($v in $values) {
        if ($v -ne $null) {
            $ws.Cells.Item($row, $col).Value = $v
        }
        $col = $col + 1
    }
}

# ---------------------------------------------------------------------
# ALC: fill in previously-empty H:N columns for rows 125-141
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ALC")

$alcRows = @{
    125 = @(905696.25, 0, 905696.25, 0, 8151266.25, $null, -8156186.25)
    126 = @(14130, 0, 14130, 0, 14130, $null, -24010)
    127 = @(35715350, 90909560, 1455.1765, 272728680, 4365.529500000001, -272723720, -14285.5295)
    128 = @(16748.75, 0, 16748.75, 0, 16748.75, $null, -26708.75)
    129 = @(1014.5574, 500.375, 1092.1698, 1501.125, 3276.5094, 3498.875, -13276.5094)
    130 = @(12000, 0, 12000, 0, 12000, $null, -22040)
    131 = @(10004262, 83332200, 4997.727, 249996600, 14993.181, -249991560, -25073.181)
    132 = @(46130.61, 52930.9, 795.3333, 158792.7, 2385.9999, -156262.7, -7445.9999)
    133 = @(38000, 0, 38000, 0, 38000, $null, -48120)
    134 = @(42000, 0, 42000, 0, 42000, $null, -52140)
    135 = @(886.2258, 643.14813, 2527, 5788.33317, 22743, -3253.33317, -27813)
    136 = @(44985, 0, 44985, 0, 44985, $null, -55185)
    137 = @(1445.2444, 1413.7307, 1488.3684, 4241.1921, 4465.1052, -1691.1921, -9565.1052)
    138 = @(18185982, 1920.037, 35720616, 5760.111, 107161848, -620.1109999999999, -107172128)
    139 = @(54475.383, 0, 54475.383, 0, 54475.383, $null, -64755.383)
    140 = @(67328.42999999999, 0, 67328.42999999999, 0, 67328.42999999999, $null, -77688.42999999999)
    141 = @(7337.1113, 2607, 13249.75, 7821, 39749.25, -2641, -50109.25)
}

foreach ($r in 125..141) {
    Set-RowValues $ws $r $alcRows[$r]
}

# ---------------------------------------------------------------------
# ARM: refresh rows 122-123
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ARM")

Set-RowValues $ws 122 @(12823038, 23811254, 3450.6667, 71433762, 10352.0001, -71431312, -15252.0001)

$ws.Cells.Item(123, 8).Value = 0
$ws.Cells.Item(123, 9).Value = 0
$ws.Cells.Item(123, 10).Value = 0
$ws.Cells.Item(123, 11).Value = 0
$ws.Cells.Item(123, 12).Value = 0
$ws.Cells.Item(123, 14).ClearContents()

# ---------------------------------------------------------------------
# BSM: refresh rows 75 and 78
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("BSM")

$ws.Cells.Item(75, 8).Value = 8516.5
$ws.Cells.Item(75, 9).Value = 2324.75
$ws.Cells.Item(75, 11).Value = 2324.75
$ws.Cells.Item(75, 13).Value = -1388.75

$ws.Cells.Item(78, 8).Value = 8516.5
$ws.Cells.Item(78, 9).Value = 2324.75
$ws.Cells.Item(78, 11).Value = 6974.25
$ws.Cells.Item(78, 13).Value = -2294.25

# ---------------------------------------------------------------------
# CUL: refresh rows 5, 122, 135
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("CUL")

Set-RowValues $ws 5   @(338165.12, 389.375, 627687.2, 1168.125, 1883061.6, -1056.125, -1883285.6)
Set-RowValues $ws 122 @(553.85297, 392.16, 1003, 3529.44, 9027, -1079.44, -13927)
Set-RowValues $ws 135 @(338165.12, 389.375, 627687.2, 3504.375, 5649184.8, -969.375, -5654254.8)

# ---------------------------------------------------------------------
# GSM: refresh row 122
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("GSM")

Set-RowValues $ws 122 @(94654.96000000001, 110275, 22802.8, 330825, 68408.39999999999, -328375, -73308.39999999999)

# ---------------------------------------------------------------------
# LTW: refresh row 46
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("LTW")

Set-RowValues $ws 46 @(63924.25, 92450.82000000001, 1165.8, 92450.82000000001, 1165.8, -92262.82000000001, -1541.8)

# ---------------------------------------------------------------------
# WVR: refresh row 136
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("WVR")

Set-RowValues $ws 136 @(7942769.5, 18523172, 7467.4585, 55569516, 22402.3755, -55566966, -27502.3755)
